$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T5").Select()

$ws.Cells.Item(3, 18).Value = $null
$ws.Cells.Item(4, 18).Value = 2021
$ws.Cells.Item(5, 18).Value = 0.9
$ws.Cells.Item(6, 18).Value = 6.5
